$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.357.12'
$ws.Range('E2').Value = '  +0.26%  '

$ws.Range('D3').Value = '3.674.04'
$ws.Range('E3').Value = '  -0.28%  '

$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').Value = '''685.29'
$ws.Range('E5').Value = '  +0.45%  '

$ws.Range('D6').Value = '''159.22'
$ws.Range('E6').Value = '  -2.12%  '

$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('D8').Value = '''0.493'
$ws.Range('E8').Value = '  -1.03%  '

$ws.Range('E9').Value = '  -1.86%  '

$ws.Range('D10').Value = '''7.04'
$ws.Range('E10').Value = '  -2.23%  '

$ws.Range('D11').Value = '''0.435'
$ws.Range('E11').Value = '  -3.55%  '

$ws.Range('E12').Value = '  -1.61%  '

$ws.Range('D13').Value = '4.292.26'
$ws.Range('E13').Value = '  -0.25%  '

$ws.Range('D14').Value = '''32.19'
$ws.Range('E14').Value = '  -4.15%  '

$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '3.665.65'
$ws.Range('E15').Value = '  -0.55%  '

$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '69.343.27'
$ws.Range('E16').Value = '  +0.14%  '

$ws.Range('E17').Value = '  +1.94%  '

$ws.Range('E18').Value = '  -3.37%  '

$ws.Range('D20').Value = '''469.43'
$ws.Range('E20').Value = '  -2.76%  '

$ws.Range('E21').Value = '  +1.03%  '

$ws.Range('D22').Value = '''0.648'
$ws.Range('E22').Value = '  -2.71%  '

$ws.Range('D23').Value = '''79.69'

$ws.Range('D24').Value = '3.819.24'
$ws.Range('E24').Value = '  -0.22%  '

$ws.Range('E25').Value = '  +0.08%  '

$ws.Range('E26').Value = '  -3.66%  '

$ws.Range('D27').Value = '''10.93'
$ws.Range('E27').Value = '  -5.51%  '

$ws.Range('D28').Value = '''9.20'
$ws.Range('E28').Value = '  -3.64%  '

$ws.Range('E29').Value = '  -1.64%  '

$ws.Range('E30').Value = '  -5.84%  '

$ws.Range('E31').Value = '  -3.07%  '

$ws.Range('D32').Value = '''1.00'
$ws.Range('E32').Value = '  +0.09%  '

$ws.Range('D33').Value = '''1.98'
$ws.Range('E33').Value = '  -5.99%  '

$ws.Range('E34').Value = '  -0.72%  '

$ws.Range('D35').Value = '3.646.03'
$ws.Range('E35').Value = '  -0.09%  '

$ws.Range('E36').Value = '  -2.34%  '

$ws.Range('D37').Value = '''8.15'
$ws.Range('E37').Value = '  -4.40%  '

$ws.Range('E38').Value = '  +1.28%  '

$ws.Range('E39').Value = '  +0.00%  '

$ws.Range('D40').Value = '''2.22'
$ws.Range('E40').Value = '  +2.50%  '

$ws.Range('D41').Value = '''0.0896'
$ws.Range('E41').Value = '  -5.28%  '

$ws.Range('E42').Value = '  -0.02%  '

$ws.Range('D43').Value = '''0.941'
$ws.Range('E43').Value = '  -1.62%  '

$ws.Range('D44').Value = '''166.30'
$ws.Range('E44').Value = '  +5.93%  '

$ws.Range('D45').Value = '''47.48'
$ws.Range('E45').Value = '  -1.22%  '

$ws.Range('D46').Value = '''0.000280'
$ws.Range('E46').Value = '  +0.87%  '

$ws.Range('D47').Value = '''2.70'
$ws.Range('E47').Value = '  -3.09%  '

$ws.Range('D48').Value = '''1.11'
$ws.Range('E48').Value = '  +4.61%  '

$ws.Range('E49').Value = '  -1.35%  '

$ws.Range('D50').Value = '''27.47'
$ws.Range('E50').Value = '  -1.63%  '

$ws.Range('E51').Value = '  -4.10%  '

Write-Host "Applied all crypto list updates"
